# Adding generated reports from pipeline run
# Insert a new "Tool Category" column (B) into the Summary sheet, shifting
# the existing severity columns (INFO/LOW/MEDIUM/HIGH/CRITICAL/UNKNOWN)
# one position to the right, and populate the new column with each tool's
# category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Shift B:G -> C:H by inserting a new column at B.
$ws.Range("B1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "Tool Category"

# Category values for each tool row.
$ws.Range("B2").Value = "SAST"
$ws.Range("B3").Value = "SCA"
$ws.Range("B4").Value = "IaC Scan"
$ws.Range("B5").Value = "Secret Scan"
